$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Spot_PT")

# Row 2 values update - automatic electricity price update
$ws.Range("A2").Value = 45946
$ws.Range("B2").Value = 117.37
$ws.Range("C2").Value = 109.67
$ws.Range("D2").Value = 106.1
$ws.Range("E2").Value = 105
$ws.Range("F2").Value = 104.98
$ws.Range("G2").Value = 106.63
$ws.Range("H2").Value = 114.89
$ws.Range("I2").Value = 128.11
$ws.Range("J2").Value = 152.89
$ws.Range("K2").Value = 124.77
$ws.Range("L2").Value = 107.52
$ws.Range("M2").Value = 89.59999999999999
$ws.Range("N2").Value = 83.09999999999999
$ws.Range("O2").Value = 81.05
$ws.Range("P2").Value = 79.41
$ws.Range("Q2").Value = 80.62
$ws.Range("R2").Value = 81.23
$ws.Range("S2").Value = 95.09999999999999
$ws.Range("T2").Value = 119.33
$ws.Range("U2").Value = 147.31
$ws.Range("V2").Value = 153.47
$ws.Range("W2").Value = 141.78
$ws.Range("X2").Value = 129.84
$ws.Range("Y2").Value = 115.22
$ws.Range("Z2").Value = 111.46
$ws.Range("AB2").Value = 135.08
$ws.Range("AD2").Value = 147.62
$ws.Range("AF2").Value = 138.83
